# Generate Report for Archive
# Updates the "Status" column value from "Ready for handoff" to "In Translation"
# across the Overview, zh-cn, and de-de worksheets, then auto-fits the
# affected columns to match the new (shorter) text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: Status for zh-cn / de-de lives in columns E and F (rows 2-3)
$overview.Range("E2:F3").Value = "In Translation"

# zh-cn sheet: Status column is C (rows 2-3)
$zhcn.Range("C2:C3").Value = "In Translation"

# de-de sheet: Status column is C (rows 2-3)
$dede.Range("C2:C3").Value = "In Translation"

# Re-fit the columns that held the status text so widths follow the new,
# shorter value (matches Excel's automatic column-width recalculation)
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5
$zhcn.Range("C1").ColumnWidth = 12.5
$dede.Range("C1").ColumnWidth = 12.5
